{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block (and\n// the blank paragraph that precedes it) that follows the\n// \"LOB1024: Mec\u00e2nica (Requisito fraco)\" paragraph near the end of the body.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"LOB1024: Mec\u00e2nica (Requisito fraco)\";\nconst footerLine1 = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst footerLine2 =\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\n\n// Find the \"LOB1024...\" paragraph, then remove the run of paragraphs that\n// follows it: the blank paragraph + the two footer paragraphs.\nlet markerIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === marker) {\n    markerIndex = i;\n    break;\n  }\n}\n\nif (markerIndex !== -1) {\n  const toDelete = [];\n  let i = markerIndex + 1;\n  // Skip/collect the single blank paragraph right after the marker.\n  if (i < paragraphs.items.length && paragraphs.items[i].text === \"\") {\n    toDelete.push(paragraphs.items[i]);\n    i++;\n  }\n  // Collect the two footer paragraphs if present.\n  if (i < paragraphs.items.length && paragraphs.items[i].text === footerLine1) {\n    toDelete.push(paragraphs.items[i]);\n    i++;\n  }\n  if (i < paragraphs.items.length && paragraphs.items[i].text === footerLine2) {\n    toDelete.push(paragraphs.items[i]);\n    i++;\n  }\n\n  for (const p of toDelete) {\n    p.delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block (and\n# the blank paragraph that precedes it) that follows the\n# \"LOB1024: Mec\u00e2nica (Requisito fraco)\" paragraph near the end of the body.\n\n$d = $word.ActiveDocument\n\n$marker = \"LOB1024: Mec\u00e2nica (Requisito fraco)\"\n$footerLine1 = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$footerLine2 = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n# Locate the \"LOB1024...\" paragraph.\n$markerIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($text -eq $marker) {\n        $markerIdx = $i\n        break\n    }\n}\n\nif ($markerIdx -ne -1) {\n    # Collect the indexes of the paragraphs to remove: the single blank\n    # paragraph right after the marker, followed by the two footer lines.\n    $toDeleteIdx = @()\n    $i = $markerIdx + 1\n\n    if ($i -le $d.Paragraphs.Count) {\n        $text = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n        if ($text -eq \"\") {\n            $toDeleteIdx += $i\n            $i++\n        }\n    }\n    if ($i -le $d.Paragraphs.Count) {\n        $text = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n        if ($text -eq $footerLine1) {\n            $toDeleteIdx += $i\n            $i++\n        }\n    }\n    if ($i -le $d.Paragraphs.Count) {\n        $text = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n        if ($text -eq $footerLine2) {\n            $toDeleteIdx += $i\n            $i++\n        }\n    }\n\n    # Delete from the last index backwards so earlier indexes stay valid.\n    for ($j = $toDeleteIdx.Count - 1; $j -ge 0; $j--) {\n        $d.Paragraphs.Item($toDeleteIdx[$j]).Range.Delete()\n    }\n}\n"}
